{"js": "// Increase font sizes throughout the resume per the commit:\n//   name header:            16pt -> 18pt  (sz 32 -> 36)\n//   everything else that is currently:\n//      9pt  -> 10pt (sz 18 -> 20)   (contact info, job dates, body/bullets)\n//     10pt  -> 11pt (sz 20 -> 22)   (Overview paragraph)\n//     11pt  -> 12pt (sz 22 -> 24)   (job titles / degree line)\n//     12pt  -> 13pt (sz 24 -> 26)   (section headers)\n//\n// Word's Office.js `Font.size` property is expressed in POINTS (not\n// half-points), so the map below is expressed in points.\nconst SIZE_MAP = {\n  16: 18,\n  9: 10,\n  10: 11,\n  11: 12,\n  12: 13,\n};\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Load text + current font size for every paragraph first.\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  p.load(\"text\");\n  p.font.load(\"size\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const text = p.text;\n  // Skip empty / whitespace-only paragraphs (nothing to resize).\n  if (!text || !text.trim()) {\n    continue;\n  }\n\n  const currentSize = Math.round(p.font.size);\n  const newSize = SIZE_MAP[currentSize];\n  if (!newSize) {\n    // Unmapped size - leave untouched.\n    continue;\n  }\n\n  // Select a range that covers exactly the paragraph's visible text\n  // (via search) rather than `paragraph.getRange()` - the latter also\n  // covers the paragraph-mark and would stamp a stray <w:rPr> onto the\n  // <w:pPr>, which is not part of the target edit. A search-based range\n  // only touches the run(s) that actually carry the text.\n  const searchResults = p.getRange().search(text, { matchCase: true });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  for (let j = 0; j < searchResults.items.length; j++) {\n    searchResults.items[j].font.size = newSize;\n  }\n  await context.sync();\n}\n", "ps1": "# Increase font sizes throughout the resume per the commit:\n#   name header:            16pt -> 18pt\n#   everything else that is currently:\n#      9pt  -> 10pt   (contact info, job dates, body/bullets)\n#     10pt  -> 11pt   (Overview paragraph)\n#     11pt  -> 12pt   (job titles / degree line)\n#     12pt  -> 13pt   (section headers)\n#\n# Word COM's Range.Font.Size is expressed in points (not half-points).\n$sizeMap = @{\n    16 = 18\n    9  = 10\n    10 = 11\n    11 = 12\n    12 = 13\n}\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $pRange = $p.Range\n    $start = $pRange.Start\n    $end = $pRange.End\n\n    # Paragraph.Range includes the trailing paragraph-mark character (\\r).\n    # Shrink the range by one character so we only touch the run(s) that\n    # actually carry visible text - setting font size on the full\n    # paragraph range (mark included) would stamp a stray run-properties\n    # element onto the paragraph mark itself, which is not part of the\n    # target edit.\n    if ($end -le $start) {\n        continue\n    }\n    $textRange = $d.Range($start, $end - 1)\n\n    if ($textRange.Text.Length -eq 0) {\n        continue\n    }\n\n    $currentSize = [int]$textRange.Font.Size\n    if ($sizeMap.ContainsKey($currentSize)) {\n        $textRange.Font.Size = $sizeMap[$currentSize]\n    }\n}\n"}
